$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.11638189569927704,
    -0.0059999999736461973,
    -0.0039999999787703189,
    -0.0079999999603170835,
    -0.0029999999805845334,
    -0.0019999999817681413,
    -0.0099999999470270495,
    0.021890850437772613,
    -0.0019999999782691624,
    -0.0019999999757800424,
    -0.0029999999711574077,
    -0.0034999999683558158,
    -0.0034999999657587821,
    -0.0079999999456852322,
    -0.00099999997532762563,
    -0.0019999999705739846,
    -0.0019999999699313875,
    -0.0039999999611506354,
    -0.0039999999822373233,
    0.0064811269900761914,
    -0.0039999999816737741,
    -0.0039999999815121257,
    -0.004999999973116509,
    -0.019999999906206156,
    -0.019999999904984911,
    -0.0024999999793564598,
    -0.0024999999790384919,
    0.075427550211641048,
    -0.0069999999527183832,
    -0.059999999722053499,
    -0.0069999999480110375,
    -0.0099999999347950563,
    -0.003999999960491607
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
